$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9705607891082764
$ws.Range("B1").Value = 1.822853207588196
$ws.Range("C1").Value = 3.173735380172729
$ws.Range("D1").Value = 3.963399887084961
$ws.Range("E1").Value = 1.123701572418213
